# Apply corrected IFRS financial figures (Kangwon Land) to company_list sheet
# Rows 2-9: numeric values replaced per corrected IFRS data; U2/U3 cells removed
# (their old "FCF" value is dropped, replaced by a single corrected value in V2/V3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("U2").ClearContents()
$ws.Range("D2").Value = 14965
$ws.Range("E2").Value = 5132
$ws.Range("F2").Value = 5132
$ws.Range("G2").Value = 4959
$ws.Range("H2").Value = 3593
$ws.Range("I2").Value = 3594
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 33751
$ws.Range("L2").Value = 6698
$ws.Range("M2").Value = 27054
$ws.Range("N2").Value = 27051
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 1070
$ws.Range("Q2").Value = 5104
$ws.Range("R2").Value = -4274
$ws.Range("S2").Value = -1480
$ws.Range("T2").Value = 802
$ws.Range("V2").Value = 4302
$ws.Range("W2").Value = 34.29
$ws.Range("X2").Value = 24.01
$ws.Range("Y2").Value = 13.76
$ws.Range("Z2").Value = 11.1
$ws.Range("AA2").Value = 24.76
$ws.Range("AB2").Value = 2590.76
$ws.Range("AC2").Value = 1680
$ws.Range("AD2").Value = 18.1
$ws.Range("AE2").Value = 13343
$ws.Range("AF2").Value = 2.28
$ws.Range("AG2").Value = 850
$ws.Range("AH2").Value = 2.8
$ws.Range("AI2").Value = 47.95
$ws.Range("AJ2").Value = 213940500

# Row 3
$ws.Range("U3").ClearContents()
$ws.Range("D3").Value = 16337
$ws.Range("E3").Value = 5954
$ws.Range("F3").Value = 5954
$ws.Range("G3").Value = 5939
$ws.Range("H3").Value = 4416
$ws.Range("I3").Value = 4416
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 36908
$ws.Range("L3").Value = 6930
$ws.Range("M3").Value = 29978
$ws.Range("N3").Value = 29975
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 1070
$ws.Range("Q3").Value = 5750
$ws.Range("R3").Value = -3989
$ws.Range("S3").Value = -1723
$ws.Range("T3").Value = 404
$ws.Range("V3").Value = 5346
$ws.Range("W3").Value = 36.45
$ws.Range("X3").Value = 27.03
$ws.Range("Y3").Value = 15.49
$ws.Range("Z3").Value = 12.5
$ws.Range("AA3").Value = 23.12
$ws.Range("AB3").Value = 2864.94
$ws.Range("AC3").Value = 2064
$ws.Range("AD3").Value = 18.6
$ws.Range("AE3").Value = 14785
$ws.Range("AF3").Value = 2.6
$ws.Range("AG3").Value = 980
$ws.Range("AH3").Value = 2.55
$ws.Range("AI3").Value = 44.99
$ws.Range("AJ3").Value = 213940500

# Row 4
$ws.Range("D4").Value = 16965
$ws.Range("E4").Value = 6186
$ws.Range("F4").Value = 6186
$ws.Range("G4").Value = 5969
$ws.Range("H4").Value = 4545
$ws.Range("I4").Value = 4546
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 39790
$ws.Range("L4").Value = 7263
$ws.Range("M4").Value = 32527
$ws.Range("N4").Value = 32525
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 1070
$ws.Range("Q4").Value = 5965
$ws.Range("R4").Value = -4115
$ws.Range("S4").Value = -1977
$ws.Range("T4").Value = 746
$ws.Range("U4").Value = 5219
$ws.Range("V4").Value = 9
$ws.Range("W4").Value = 36.46
$ws.Range("X4").Value = 26.79
$ws.Range("Y4").Value = 14.55
$ws.Range("Z4").Value = 11.85
$ws.Range("AA4").Value = 22.33
$ws.Range("AB4").Value = 3103.31
$ws.Range("AC4").Value = 2125
$ws.Range("AD4").Value = 16.82
$ws.Range("AE4").Value = 16043
$ws.Range("AF4").Value = 2.23
$ws.Range("AG4").Value = 990
$ws.Range("AH4").Value = 2.77
$ws.Range("AI4").Value = 44.15
$ws.Range("AJ4").Value = 213940500

# Row 5
$ws.Range("D5").Value = 15478
$ws.Range("E5").Value = 5309
$ws.Range("F5").Value = 5309
$ws.Range("G5").Value = 5662
$ws.Range("H5").Value = 4375
$ws.Range("I5").Value = 4376
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 41857
$ws.Range("L5").Value = 6844
$ws.Range("M5").Value = 35012
$ws.Range("N5").Value = 35011
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1070
$ws.Range("Q5").Value = 4630
$ws.Range("R5").Value = -2502
$ws.Range("S5").Value = -2003
$ws.Range("T5").Value = 855
$ws.Range("U5").Value = 3776
$ws.Range("V5").Value = 14
$ws.Range("W5").Value = 34.3
$ws.Range("X5").Value = 28.27
$ws.Range("Y5").Value = 12.96
$ws.Range("Z5").Value = 10.72
$ws.Range("AA5").Value = 19.55
$ws.Range("AB5").Value = 3335.69
$ws.Range("AC5").Value = 2046
$ws.Range("AD5").Value = 17.01
$ws.Range("AE5").Value = 17269
$ws.Range("AF5").Value = 2.02
$ws.Range("AG5").Value = 990
$ws.Range("AH5").Value = 2.84
$ws.Range("AI5").Value = 45.86
$ws.Range("AJ5").Value = 213940500

# Row 6
$ws.Range("D6").Value = 14381
$ws.Range("E6").Value = 4307
$ws.Range("F6").Value = 4307
$ws.Range("G6").Value = 4107
$ws.Range("H6").Value = 2972
$ws.Range("I6").Value = 2973
$ws.Range("K6").Value = 42312
$ws.Range("L6").Value = 6445
$ws.Range("M6").Value = 35868
$ws.Range("N6").Value = 35867
$ws.Range("P6").Value = 1070
$ws.Range("Q6").Value = 3639
$ws.Range("R6").Value = -2020
$ws.Range("S6").Value = -2024
$ws.Range("T6").Value = 1033
$ws.Range("U6").Value = 2607
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 29.95
$ws.Range("X6").Value = 20.67
$ws.Range("Y6").Value = 8.390000000000001
$ws.Range("Z6").Value = 7.06
$ws.Range("AA6").Value = 17.97
$ws.Range("AB6").Value = 3415.78
$ws.Range("AC6").Value = 1390
$ws.Range("AD6").Value = 23.03
$ws.Range("AE6").Value = 17691
$ws.Range("AF6").Value = 1.81
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 2.81
$ws.Range("AI6").Value = 61.37
$ws.Range("AJ6").Value = 213940500

# Row 7
$ws.Range("D7").Value = 15238
$ws.Range("E7").Value = 5172
$ws.Range("G7").Value = 4939
$ws.Range("H7").Value = 3407
$ws.Range("I7").Value = 3408
$ws.Range("K7").Value = 43701
$ws.Range("L7").Value = 6422
$ws.Range("M7").Value = 37279
$ws.Range("N7").Value = 37279
$ws.Range("P7").Value = 1070
$ws.Range("Q7").Value = 3934
$ws.Range("R7").Value = 962
$ws.Range("S7").Value = -1834
$ws.Range("T7").Value = 580
$ws.Range("U7").Value = 3387
$ws.Range("W7").Value = 33.94
$ws.Range("X7").Value = 22.36
$ws.Range("Y7").Value = 9.32
$ws.Range("Z7").Value = 7.92
$ws.Range("AA7").Value = 17.23
$ws.Range("AC7").Value = 1593
$ws.Range("AD7").Value = 17.3
$ws.Range("AE7").Value = 18388
$ws.Range("AF7").Value = 1.5
$ws.Range("AG7").Value = 918
$ws.Range("AH7").Value = 3.33
$ws.Range("AI7").Value = 57.64

# Row 8
$ws.Range("D8").Value = 15923
$ws.Range("E8").Value = 5227
$ws.Range("G8").Value = 5393
$ws.Range("H8").Value = 3994
$ws.Range("I8").Value = 3995
$ws.Range("K8").Value = 45988
$ws.Range("L8").Value = 6598
$ws.Range("M8").Value = 39390
$ws.Range("N8").Value = 39391
$ws.Range("P8").Value = 1070
$ws.Range("Q8").Value = 4762
$ws.Range("R8").Value = -1364
$ws.Range("S8").Value = -1877
$ws.Range("T8").Value = 409
$ws.Range("U8").Value = 4058
$ws.Range("W8").Value = 32.82
$ws.Range("X8").Value = 25.09
$ws.Range("Y8").Value = 10.42
$ws.Range("Z8").Value = 8.91
$ws.Range("AA8").Value = 16.75
$ws.Range("AC8").Value = 1867
$ws.Range("AD8").Value = 14.75
$ws.Range("AE8").Value = 19429
$ws.Range("AF8").Value = 1.42
$ws.Range("AG8").Value = 955
$ws.Range("AH8").Value = 3.47
$ws.Range("AI8").Value = 51.12

# Row 9
$ws.Range("D9").Value = 16791
$ws.Range("E9").Value = 5542
$ws.Range("G9").Value = 5983
$ws.Range("H9").Value = 4273
$ws.Range("I9").Value = 4274
$ws.Range("K9").Value = 48470
$ws.Range("L9").Value = 6772
$ws.Range("M9").Value = 41698
$ws.Range("N9").Value = 41699
$ws.Range("P9").Value = 1070
$ws.Range("Q9").Value = 5047
$ws.Range("R9").Value = -1152
$ws.Range("S9").Value = -1932
$ws.Range("T9").Value = 373
$ws.Range("U9").Value = 4418
$ws.Range("W9").Value = 33.01
$ws.Range("X9").Value = 25.45
$ws.Range("Y9").Value = 10.54
$ws.Range("Z9").Value = 9.050000000000001
$ws.Range("AA9").Value = 16.24
$ws.Range("AC9").Value = 1998
$ws.Range("AD9").Value = 13.79
$ws.Range("AE9").Value = 20568
$ws.Range("AF9").Value = 1.34
$ws.Range("AG9").Value = 1001
$ws.Range("AH9").Value = 3.63
$ws.Range("AI9").Value = 50.1

Write-Host "Updated IFRS rows 2-9 on company_list sheet"
